$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2438
$ws.Range("I2").Value = 2297.6
$ws.Range("J2").Value = 2578.4
$ws.Range("K2").Value = 2297.6
$ws.Range("L2").Value = 2578.4
$ws.Range("M2").Value = -2184.6
$ws.Range("N2").Value = -2804.4

$ws.Range("H12").Value = 238.71428
$ws.Range("I12").Value = 174.66667
$ws.Range("K12").Value = 174.66667
$ws.Range("M12").Value = -4.666670000000011

$ws.Range("H29").Value = 3980
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15562

$ws.Range("H41").Value = 294.77777
$ws.Range("I41").Value = 236
$ws.Range("J41").Value = 500.5
$ws.Range("K41").Value = 236
$ws.Range("L41").Value = 500.5
$ws.Range("M41").Value = 204
$ws.Range("N41").Value = -1380.5

$ws.Range("H43").Value = 1686.0605
$ws.Range("I43").Value = 2079
$ws.Range("J43").Value = 1660.7097
$ws.Range("K43").Value = 2079
$ws.Range("L43").Value = 1660.7097
$ws.Range("M43").Value = -2010
$ws.Range("N43").Value = -1798.7097

$ws.Range("H58").Value = 543.8
$ws.Range("I58").Value = 543.8
$ws.Range("K58").Value = 1631.4
$ws.Range("M58").Value = -1481.4

$ws.Range("H76").Value = 774377.4
$ws.Range("J76").Value = 7209.8
$ws.Range("L76").Value = 7209.8
$ws.Range("N76").Value = -7839.8

$ws.Range("H79").Value = 774377.4
$ws.Range("J79").Value = 7209.8
$ws.Range("L79").Value = 7209.8
$ws.Range("N79").Value = -9393.799999999999

$ws.Range("H82").Value = 2291.6
$ws.Range("I82").Value = 352.75
$ws.Range("K82").Value = 1058.25
$ws.Range("M82").Value = -652.25

$ws.Range("H85").Value = 2291.6
$ws.Range("I85").Value = 352.75
$ws.Range("K85").Value = 1058.25
$ws.Range("M85").Value = 345.75

$ws.Range("H86").Value = 772980
$ws.Range("J86").Value = 6264.3335
$ws.Range("L86").Value = 6264.3335
$ws.Range("N86").Value = -8510.333500000001

$ws.Range("H88").Value = 3629.6667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3629.6667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3629.6667
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4441.6667

$ws.Range("H89").Value = 772980
$ws.Range("J89").Value = 6264.3335
$ws.Range("L89").Value = 31321.6675
$ws.Range("N89").Value = -42553.6675

$ws.Range("H91").Value = 3629.6667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3629.6667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3629.6667
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6437.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7171.7144
$ws.Range("J46").Value = 7774.5
$ws.Range("L46").Value = 7774.5
$ws.Range("N46").Value = -8412.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8358.893
$ws.Range("I86").Value = 6176.5
$ws.Range("J86").Value = 11268.75
$ws.Range("K86").Value = 6176.5
$ws.Range("L86").Value = 11268.75
$ws.Range("M86").Value = -5053.5
$ws.Range("N86").Value = -13514.75

$ws.Range("H89").Value = 8358.893
$ws.Range("I89").Value = 6176.5
$ws.Range("J89").Value = 11268.75
$ws.Range("K89").Value = 30882.5
$ws.Range("L89").Value = 56343.75
$ws.Range("M89").Value = -25266.5
$ws.Range("N89").Value = -67575.75

$ws.Range("H134").Value = 6555.875
$ws.Range("I134").Value = 4767.5557
$ws.Range("K134").Value = 14302.6671
$ws.Range("M134").Value = -11767.6671

$ws.Range("H141").Value = 78798.60000000001
$ws.Range("J141").Value = 78798.60000000001
$ws.Range("L141").Value = 78798.60000000001
$ws.Range("N141").Value = -89158.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2848.1365
$ws.Range("J31").Value = 4570.7
$ws.Range("L31").Value = 4570.7
$ws.Range("N31").Value = -5160.7

$ws.Range("H34").Value = 2848.1365
$ws.Range("J34").Value = 4570.7
$ws.Range("L34").Value = 4570.7
$ws.Range("N34").Value = -4974.7

$ws.Range("H62").Value = 6544.4
$ws.Range("I62").Value = 2222
$ws.Range("J62").Value = 7625
$ws.Range("K62").Value = 2222
$ws.Range("L62").Value = 7625
$ws.Range("M62").Value = -1598
$ws.Range("N62").Value = -8873

$ws.Range("H65").Value = 6544.4
$ws.Range("I65").Value = 2222
$ws.Range("J65").Value = 7625
$ws.Range("K65").Value = 11110
$ws.Range("L65").Value = 38125
$ws.Range("M65").Value = -7990
$ws.Range("N65").Value = -44365

$ws.Range("H122").Value = 2681.6875
$ws.Range("I122").Value = 1940
$ws.Range("J122").Value = 4906.75
$ws.Range("K122").Value = 5820
$ws.Range("L122").Value = 14720.25
$ws.Range("M122").Value = -3370
$ws.Range("N122").Value = -19620.25

$ws.Range("H139").Value = 111000
$ws.Range("J139").Value = 111000
$ws.Range("L139").Value = 111000
$ws.Range("N139").Value = -121280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 74
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 74
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 854.44446
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 854.44446
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -559.44446
$ws.Range("N22").Value = -3590

$ws.Range("I27").Value = 854.44446
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 854.44446
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -747.44446
$ws.Range("N27").Value = -3214

$ws.Range("H40").Value = 6494.923
$ws.Range("I40").Value = 4354.875
$ws.Range("K40").Value = 4354.875
$ws.Range("M40").Value = -4218.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28497.666
$ws.Range("J62").Value = 42750
$ws.Range("L62").Value = 42750
$ws.Range("N62").Value = -43998

$ws.Range("H65").Value = 28497.666
$ws.Range("J65").Value = 42750
$ws.Range("L65").Value = 213750
$ws.Range("N65").Value = -219990

$ws.Range("H122").Value = 4009.7837
$ws.Range("I122").Value = 3909.5356
$ws.Range("J122").Value = 4321.6665
$ws.Range("K122").Value = 11728.6068
$ws.Range("L122").Value = 12964.9995
$ws.Range("M122").Value = -9278.606800000001
$ws.Range("N122").Value = -17864.9995

$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
